$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text (matching source data),
# since Excel would otherwise auto-convert a clean single-decimal string to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cryptocurrency price and volume figures.
$ws.Range("D2").Value = '43.797.21'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '2.294.17'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '97.59'
$ws.Range("E5").Value = '  +2.94%  '
$ws.Range("D6").Value = '270.00'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("D10").Value = '45.34'
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("E13").Value = '  +1.97%  '
$ws.Range("D14").Value = '15.84'
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '2.638.16'
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").Value = '0.859'
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '2.292.73'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").Value = '43.798.46'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").Value = '72.24'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("E22").Value = '  +9.46%  '
$ws.Range("D23").Value = '232.03'
$ws.Range("E23").Value = '  -3.57%  '
$ws.Range("D24").Value = '9.14'
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("E25").Value = '  +9.15%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = '38.24'
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("D31").Value = '176.59'
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("D32").Value = '21.84'
$ws.Range("E32").Value = '  -2.68%  '
$ws.Range("D33").Value = '0.0893'
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("E35").Value = '  +1.38%  '
$ws.Range("D36").Value = '4.70'
$ws.Range("E36").Value = '  +7.55%  '
$ws.Range("E37").Value = '  +2.38%  '
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("D39").Value = '3.50'
$ws.Range("E39").Value = '  +4.43%  '
$ws.Range("D40").Value = '0.236'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").Value = '1.37'
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("D43").Value = '12.18'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").Value = '64.77'
$ws.Range("E44").Value = '  +5.07%  '
$ws.Range("D45").Value = '5.23'
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("D46").Value = '8.77'
$ws.Range("E46").Value = '  -4.28%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("D49").Value = '99.16'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("E51").Value = '  +11.07%  '
